$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet to "Products"
$ws.Name = "Products"

# Correct the header row (the stored values in A1/B1 were "a"/"b", not the
# unused "Item"/"URL" shared strings) and add the new "Price" header
$ws.Range("A1").Value = "Item"
$ws.Range("B1").Value = "URL"
$ws.Range("C1").Value = "Price"

# Add new product row
$ws.Range("A2").Value = "Nike Force 1"
$ws.Range("B2").Value = "https://www.hepsiburada.com/nike-erkek-yuruyus-ayakkabisi-force-1-cw2288-p-HBCV00001CN5ZW?magaza=NSPORT"
